# Scheduled-runner market data refresh for Pandaemonium_Profits.
# Re-pulls currentAveragePrice(NQ/HQ) + LevePrice(NQ/HQ) + LeveProfit(NQ/HQ)
# columns (H:N) per leve row on each Disciple of the Hand sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 4
$ws.Range("H4").Value = 2750.625
$ws.Range("I4").Value = 3087.8572
$ws.Range("J4").Value = 390
$ws.Range("K4").Value = 3087.8572
$ws.Range("L4").Value = 390
$ws.Range("M4").Value = -2973.8572
$ws.Range("N4").Value = -618
# row 18
$ws.Range("H18").Value = 329.33334
$ws.Range("I18").Value = 329.33334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 329.33334
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -45.33334000000002
$ws.Range("N18").ClearContents()
# row 53
$ws.Range("H53").Value = 149.2069
$ws.Range("I53").Value = 116.30769
$ws.Range("K53").Value = 116.30769
$ws.Range("M53").Value = 520.69231
# row 92
$ws.Range("H92").Value = 536.8
$ws.Range("I92").Value = 536.8
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 536.8
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 711.2
$ws.Range("N92").ClearContents()
# row 107
$ws.Range("H107").Value = 1145.4
$ws.Range("I107").Value = 1028.8462
$ws.Range("J107").Value = 1903
$ws.Range("K107").Value = 1028.8462
$ws.Range("L107").Value = 1903
$ws.Range("M107").Value = 891.1538
$ws.Range("N107").Value = -5743
# row 115
$ws.Range("H115").Value = 1878.5
$ws.Range("I115").Value = 1433.5714
$ws.Range("J115").Value = 2916.6667
$ws.Range("K115").Value = 4300.7142
$ws.Range("L115").Value = 8750.000100000001
$ws.Range("M115").Value = -2733.7142
$ws.Range("N115").Value = -11884.0001
# row 129
$ws.Range("H129").Value = 1099.6389
$ws.Range("J129").Value = 1132.9565
$ws.Range("L129").Value = 3398.8695
$ws.Range("N129").Value = -13398.8695
# row 141
$ws.Range("H141").Value = 1909.7646
$ws.Range("I141").Value = 1407.1464
$ws.Range("J141").Value = 3970.5
$ws.Range("K141").Value = 4221.439200000001
$ws.Range("L141").Value = 11911.5
$ws.Range("M141").Value = 958.5607999999993
$ws.Range("N141").Value = -22271.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 18554.451
$ws.Range("I32").Value = 20393.861
$ws.Range("J32").Value = 11442.066
$ws.Range("K32").Value = 20393.861
$ws.Range("L32").Value = 11442.066
$ws.Range("M32").Value = -20106.861
$ws.Range("N32").Value = -12016.066
# row 43
$ws.Range("H43").Value = 15377
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15377
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15377
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -16003
# row 102
$ws.Range("H102").Value = 1197459.9
$ws.Range("I102").Value = 1854107.2
$ws.Range("J102").Value = 3555.5454
$ws.Range("K102").Value = 1854107.2
$ws.Range("L102").Value = 3555.5454
$ws.Range("M102").Value = -1852485.2
$ws.Range("N102").Value = -6799.5454
# row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 1162
$ws.Range("I99").Value = 1133.3334
$ws.Range("J99").Value = 1205
$ws.Range("K99").Value = 1133.3334
$ws.Range("L99").Value = 1205
$ws.Range("M99").Value = 364.6666
$ws.Range("N99").Value = -4201
# row 132
$ws.Range("H132").Value = 49015.5
$ws.Range("J132").Value = 61450.715
$ws.Range("L132").Value = 61450.715
$ws.Range("N132").Value = -71570.715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 3250651.2
$ws.Range("I58").Value = 6996507
$ws.Range("J58").Value = 4242.8
$ws.Range("K58").Value = 6996507
$ws.Range("L58").Value = 4242.8
$ws.Range("M58").Value = -6996304
$ws.Range("N58").Value = -4648.8
# row 136
$ws.Range("H136").Value = 3250651.2
$ws.Range("I136").Value = 6996507
$ws.Range("J136").Value = 4242.8
$ws.Range("K136").Value = 20989521
$ws.Range("L136").Value = 12728.4
$ws.Range("M136").Value = -20986971
$ws.Range("N136").Value = -17828.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 122
$ws.Range("H122").Value = 916
$ws.Range("I122").Value = 463.66666
$ws.Range("J122").Value = 1096.9333
$ws.Range("K122").Value = 4172.99994
$ws.Range("L122").Value = 9872.3997
$ws.Range("M122").Value = -1722.99994
$ws.Range("N122").Value = -14772.3997
# row 131
$ws.Range("H131").Value = 1055.322
$ws.Range("I131").Value = 931.25
$ws.Range("J131").Value = 1074.7843
$ws.Range("K131").Value = 2793.75
$ws.Range("L131").Value = 3224.3529
$ws.Range("M131").Value = 2246.25
$ws.Range("N131").Value = -13304.3529
# row 132
$ws.Range("H132").Value = 1994.6111
$ws.Range("I132").Value = 3030.8
$ws.Range("J132").Value = 1596.0769
$ws.Range("K132").Value = 27277.2
$ws.Range("L132").Value = 14364.6921
$ws.Range("M132").Value = -24747.2
$ws.Range("N132").Value = -19424.6921

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 4617.0347
$ws.Range("I7").Value = 4499.2
$ws.Range("J7").Value = 4878.8887
$ws.Range("K7").Value = 4499.2
$ws.Range("L7").Value = 4878.8887
$ws.Range("M7").Value = -4387.2
$ws.Range("N7").Value = -5102.8887
# row 16
$ws.Range("H16").Value = 938.2273
$ws.Range("I16").Value = 656.41174
$ws.Range("J16").Value = 1896.4
$ws.Range("K16").Value = 656.41174
$ws.Range("L16").Value = 1896.4
$ws.Range("M16").Value = -486.41174
$ws.Range("N16").Value = -2236.4
# row 22
$ws.Range("H22").Value = 991.64703
$ws.Range("I22").Value = 1223.5
$ws.Range("J22").Value = 920.3077
$ws.Range("K22").Value = 1223.5
$ws.Range("L22").Value = 920.3077
$ws.Range("M22").Value = -928.5
$ws.Range("N22").Value = -1510.3077
# row 27
$ws.Range("H27").Value = 991.64703
$ws.Range("I27").Value = 1223.5
$ws.Range("J27").Value = 920.3077
$ws.Range("K27").Value = 1223.5
$ws.Range("L27").Value = 920.3077
$ws.Range("M27").Value = -1116.5
$ws.Range("N27").Value = -1134.3077
# row 55
$ws.Range("H55").Value = 476.9375
$ws.Range("I55").Value = 448.77777
$ws.Range("J55").Value = 513.1429000000001
$ws.Range("K55").Value = 448.77777
$ws.Range("L55").Value = 513.1429000000001
$ws.Range("M55").Value = -275.77777
$ws.Range("N55").Value = -859.1429000000001
# row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# row 82
$ws.Range("H82").Value = 2373.6
$ws.Range("I82").Value = 2251.6
$ws.Range("J82").Value = 2739.6
$ws.Range("K82").Value = 2251.6
$ws.Range("L82").Value = 2739.6
$ws.Range("M82").Value = -1890.6
$ws.Range("N82").Value = -3461.6
# row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# row 85
$ws.Range("H85").Value = 2373.6
$ws.Range("I85").Value = 2251.6
$ws.Range("J85").Value = 2739.6
$ws.Range("K85").Value = 2251.6
$ws.Range("L85").Value = 2739.6
$ws.Range("M85").Value = -1003.6
$ws.Range("N85").Value = -5235.6
# row 93
$ws.Range("H93").Value = 657
$ws.Range("I93").Value = 608
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 608
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 640
$ws.Range("N93").Value = -3496
# row 126
$ws.Range("H126").Value = 4617.0347
$ws.Range("I126").Value = 4499.2
$ws.Range("J126").Value = 4878.8887
$ws.Range("K126").Value = 13497.6
$ws.Range("L126").Value = 14636.6661
$ws.Range("M126").Value = -11027.6
$ws.Range("N126").Value = -19576.6661
# row 136
$ws.Range("H136").Value = 3300.3713
$ws.Range("I136").Value = 1744.591
$ws.Range("J136").Value = 5933.231
$ws.Range("K136").Value = 5233.772999999999
$ws.Range("L136").Value = 17799.693
$ws.Range("M136").Value = -2683.772999999999
$ws.Range("N136").Value = -22899.693

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 3527
$ws.Range("I62").Value = 3750
$ws.Range("J62").Value = 3458.3845
$ws.Range("K62").Value = 3750
$ws.Range("L62").Value = 3458.3845
$ws.Range("M62").Value = -3126
$ws.Range("N62").Value = -4706.3845
# row 65
$ws.Range("H65").Value = 3527
$ws.Range("I65").Value = 3750
$ws.Range("J65").Value = 3458.3845
$ws.Range("K65").Value = 18750
$ws.Range("L65").Value = 17291.9225
$ws.Range("M65").Value = -15630
$ws.Range("N65").Value = -23531.9225
# row 132
$ws.Range("H132").Value = 2073.543
$ws.Range("I132").Value = 813.3
$ws.Range("J132").Value = 3753.8667
$ws.Range("K132").Value = 2439.9
$ws.Range("L132").Value = 11261.6001
$ws.Range("M132").Value = 90.10000000000036
$ws.Range("N132").Value = -16321.6001
# row 136
$ws.Range("H136").Value = 4261.3877
$ws.Range("I136").Value = 3193.7346
$ws.Range("J136").Value = 5948.968
$ws.Range("K136").Value = 9581.203799999999
$ws.Range("L136").Value = 17846.904
$ws.Range("M136").Value = -7031.203799999999
$ws.Range("N136").Value = -22946.904
